$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.915.44"
$ws.Range("E2").Value = "  +1.42%  "

$ws.Range("D3").Value = "2.361.20"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("E5").Value = "  +5.93%  "

$ws.Range("D6").Value = "'241.52"
$ws.Range("E6").Value = "  +2.94%  "

$ws.Range("D7").Value = "'76.45"
$ws.Range("E7").Value = "  +5.29%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").Value = "'0.636"
$ws.Range("E9").Value = "  +27.90%  "

$ws.Range("E10").Value = "  +5.23%  "

$ws.Range("D11").Value = "'57.34"
$ws.Range("E11").Value = "  +0.83%  "

$ws.Range("D12").Value = "'33.44"
$ws.Range("E12").Value = "  +22.32%  "

$ws.Range("D13").Value = "'7.59"
$ws.Range("E13").Value = "  +20.36%  "

$ws.Range("E14").Value = "  +2.00%  "

$ws.Range("D15").Value = "2.712.26"
$ws.Range("E15").Value = "  +0.07%  "

$ws.Range("D16").Value = "'16.93"
$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("E17").Value = "  +6.79%  "

$ws.Range("D18").Value = "2.363.07"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").Value = "43.919.02"
$ws.Range("E19").Value = "  +1.40%  "

$ws.Range("E20").Value = "  +4.46%  "

$ws.Range("D21").Value = "'6.69"
$ws.Range("E21").Value = "  +4.75%  "

$ws.Range("D22").Value = "'77.59"
$ws.Range("E22").Value = "  +3.07%  "

$ws.Range("D23").Value = "'257.06"
$ws.Range("E23").Value = "  +2.41%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("E25").Value = "  +2.26%  "

$ws.Range("D26").Value = "'11.20"
$ws.Range("E26").Value = "  +11.17%  "

$ws.Range("E27").Value = "  -6.89%  "

$ws.Range("E28").Value = "  +15.22%  "

$ws.Range("E29").Value = "  +1.39%  "

$ws.Range("D30").Value = "'23.18"
$ws.Range("E30").Value = "  +2.59%  "

$ws.Range("D31").Value = "'174.99"
$ws.Range("E31").Value = "  +1.34%  "

$ws.Range("E32").Value = "  -2.06%  "

$ws.Range("E33").Value = "  +5.72%  "

$ws.Range("D34").Value = "'5.33"
$ws.Range("E34").Value = "  +6.35%  "

$ws.Range("E35").Value = "  +8.50%  "

$ws.Range("E36").Value = "  +6.25%  "

$ws.Range("D37").Value = "'3.83"
$ws.Range("E37").Value = "  +1.95%  "

$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("D39").Value = "'6.48"
$ws.Range("E39").Value = "  -0.87%  "

$ws.Range("D40").Value = "'0.0277"
$ws.Range("E40").Value = "  +8.06%  "

$ws.Range("D41").Value = "'19.26"
$ws.Range("E41").Value = "  -0.47%  "

$ws.Range("E42").Value = "  +18.27%  "

$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("D44").Value = "'8.94"
$ws.Range("E44").Value = "  -0.22%  "

$ws.Range("E45").Value = "  +5.95%  "

$ws.Range("E46").Value = "  +13.49%  "

$ws.Range("E47").Value = "  +5.11%  "

$ws.Range("E48").Value = "  +2.60%  "

$ws.Range("D49").Value = "'102.09"
$ws.Range("E49").Value = "  +2.44%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "'54.97"
$ws.Range("E51").Value = "  +8.39%  "
